$d = $word.ActiveDocument

# Locate the unique sentence we need to edit.
$full = $d.Content
$full.Find.Execute(
    "We will then use an RNA Protein detection assay such as the ",
    $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
$start = $full.Start

# Remove "then " (offset 8..13 relative to the start of the sentence):
# "We will then use an RNA Protein..." -> "We will use an RNA Protein..."
$rThen = $d.Range($start + 8, $start + 13)
$rThen.Text = ""

# After the deletion, the "P" of "Protein" sits 19 characters after $start.
# Change "Protein" -> "protein", forcing the run to split into three runs
# (identical formatting, matching how Word splits a run when its text is
# edited mid-run) by nudging a character-formatting property before and
# after the text edit.
$rP = $d.Range($start + 19, $start + 20)
$rP.Text = "p"
$rP.Font.Size = 20
$rP.Font.Size = 9
